$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the current row 75,
# pushing the existing rows 75-107 down to 76-108.
$ws.Rows.Item(75).Insert()

$ws.Cells.Item(75, 1).Value = 5
$ws.Cells.Item(75, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(75, 3).Value = "Maule"
$ws.Cells.Item(75, 4).Value = 44455
$ws.Cells.Item(75, 5).Value = 7
$ws.Cells.Item(75, 6).Value = 100112017
$ws.Cells.Item(75, 7).Value = "Apio"
$ws.Cells.Item(75, 8).Value = "Americana (o)"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 300
$ws.Cells.Item(75, 11).Value = 7500
$ws.Cells.Item(75, 12).Value = 7500
$ws.Cells.Item(75, 13).Value = 7500
$ws.Cells.Item(75, 14).Value = "$/docena de matas"
$ws.Cells.Item(75, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(75, 16).Value = 1250
$ws.Cells.Item(75, 17).Value = 6
$ws.Cells.Item(75, 18).Value = "Hortaliza"
